$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert 3 new rows before the old row 248 (gui/menu/research/description/well_contruction_lvl_3) ---
$ws.Rows("248:250").Insert()

$ws.Range("A248").Value = "gui/menu/research/description/towers_artillery_acid"
$ws.Range("B248").Value = "Artillery utilizing acidic rounds"

$ws.Range("A249").Value = "gui/menu/research/description/towers_artillery_napalm"
$ws.Range("B249").Value = "Artillery using indicidary rounds"

$ws.Range("A250").Value = "gui/menu/research/description/towers_rocket"
$ws.Range("B250").Value = "Simple rocket powered granade defense towers"

# --- Insert 2 new rows before old row 314 (now shifted to row 317: gui/menu/research/name/towers_lowcaliber) ---
$ws.Rows("317:318").Insert()

$ws.Range("A317").Value = "gui/menu/research/name/towers_artillery_acid"
$ws.Range("B317").Value = "Acidic Artillery"

$ws.Range("A318").Value = "gui/menu/research/name/towers_artillery_napalm"
$ws.Range("B318").Value = "Incidiary Artillery"

# --- Insert 3 new rows after that same row (now row 319), i.e. before old row 315 (now shifted to row 320) ---
$ws.Rows("320:322").Insert()

$ws.Range("A320").Value = "gui/menu/research/name/towers_rocket"
$ws.Range("B320").Value = "Rocket Towers"

$ws.Range("A321").Value = "gui/menu/research/name/towers_rocket_lvl_2"
$ws.Range("B321").Value = "Rocket Towers 2"

$ws.Range("A322").Value = "gui/menu/research/name/towers_rocket_lvl_3"
$ws.Range("B322").Value = "Rocket Towers 3"

# --- Update the active selection to match the final edit location ---
$ws.Range("A251").Select()
